$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 4 ("Supervised Learning: Overview"): merge the "From data points,
# predict continuous valued " / "outputs." runs in the Regression bullet
# into a single run with unified text.
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$tr4 = $sh4.TextFrame.TextRange
$full4 = $tr4.Text
$needle4 = "From data points, predict continuous valued outputs."
$idx4 = $full4.IndexOf($needle4)
$sub4 = $tr4.Characters($idx4 + 1, $needle4.Length)
$sub4.Text = $needle4

# ---------------------------------------------------------------------------
# Slide 5 ("Supervised Learning"): TextBox 10 grows taller (its text now
# wraps to an extra line) and its final line gains "/covariates".
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(9)

$tr5 = $sh5.TextFrame.TextRange
$full5 = $tr5.Text
$needle5 = "independent variables"
$idx5 = $full5.IndexOf($needle5)
$sub5 = $tr5.Characters($idx5 + 1, $needle5.Length)
$sub5.Text = "independent variables/covariates"

# Resize the textbox: cy 1446550 -> 1785104 EMU (cx unchanged). Must happen
# after the text edit above, since this shape auto-fits and an explicit
# height set before a text edit gets clobbered by the subsequent reflow.
$sh5.Height = 140.55941772460938

# ---------------------------------------------------------------------------
# Slide 8 ("Least Squares with Multiple Inputs"): TextBox 11 is repositioned.
# off x 4596549 -> 4624258, off y 1678426 -> 1380387 EMU (ext unchanged).
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(9)
$sh8.Left = 364.1148376464844
$sh8.Top = 108.69193267822266
